$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before column B. This shifts the old "B" (name)
# and "C" (major-code) columns to "D" and "E", and shifts the old "F"/"G"
# lookup-table columns to "H"/"I" - matching the structural shift seen in
# the diff.
$ws.Columns("B:C").Insert()

# Split what used to be a single "name" column into three columns:
# prefix / first name / last name.
$ws.Range("B1").Value = "คำนำหน้า"
$ws.Range("C1").Value = "ชื่อ"
$ws.Range("D1").Value = "นามสุกล"

# Give the newly inserted columns (B, C) the same custom width as column A
# (they line up with column A visually, even though this engine cannot
# reproduce the exact fractional width Excel stores for column A).
$ws.Columns("B:C").ColumnWidth = $ws.Columns("A").ColumnWidth()

# Update the selection / active cell to match the new layout.
$null = $ws.Range("A2:E11").Select()

# Configure the printed page (A4, portrait) which adds a <pageSetup> node.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
